$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.704.54'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.600.39'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '211.39'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").Value = '1.825.08'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '1.651.41'
$ws.Range("E13").Value = '  +4.62%  '
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '65.32'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '26.683.05'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("E18").Value = '  +4.33%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '209.47'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '7.17'
$ws.Range("E21").Value = '  +3.14%  '
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").Value = '2.29'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").Value = '142.87'
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").Value = '15.36'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '0.0517'
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").Value = '1.293.37'
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("E35").Value = '  -4.78%  '
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  +18.92%  '
$ws.Range("D40").Value = '0.826'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = '63.21'
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").Value = '1.736.87'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").Value = '91.09'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("E51").Value = '  +0.23%  '
